# membership-card.xlsx: correct the membershipId sample value in row 2.
# The cell used to hold the number 1; it should hold the text "1a"
# (membership ids are alphanumeric, not purely numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "1a"

# Leave the cursor where the author left it when the file was saved.
[void]$ws.Range("C3").Select()
